# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", populated
#   with the per-fund holdings for that quarter.
# - Update the "总计" (totals) sheet: insert a new top data row for
#   "2022-Q1" (4 funds held, 0.28 亿元 market value), pushing the existing
#   2021-Q4 / 2021-Q3 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" worksheet (content + formatting) while it is
#    still under its default temporary name/position. We rename + move it
#    to its final place only once it is fully populated, because sheet
#    references here are position-bound and are invalidated by Move().
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# Clone the header-row (B1:H1) and the "index" column (A) formatting from
# the existing "2021-Q4" sheet so the new sheet matches the workbook's
# established look (bold, centered, thin-bordered header; bold bordered
# index column).
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# The fund-code / numeric-looking text columns must stay TEXT (leading
# zeros in codes, fixed-decimal display values) rather than being coerced
# to numbers, so force a text number format before assigning them.
$newSheet.Range("B2:G5").NumberFormat = "@"

# Row 2 - 012473 / 大成成长回报六个月持有期混合型证券投资基金A
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "012473"
$newSheet.Range("C2").Value = "大成成长回报六个月持有期混合型证券投资基金A"
$newSheet.Range("D2").Value = "8.97"
$newSheet.Range("E2").Value = "71.30"
$newSheet.Range("F2").Value = "2.49"
$newSheet.Range("G2").Value = "0.2234"
$newSheet.Range("H2").Value = 10

# Row 3 - 011073 / 鹏华安润混合A
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "011073"
$newSheet.Range("C3").Value = "鹏华安润混合A"
$newSheet.Range("D3").Value = "3.07"
$newSheet.Range("E3").Value = "29.79"
$newSheet.Range("F3").Value = "1.28"
$newSheet.Range("G3").Value = "0.0393"
$newSheet.Range("H3").Value = 3

# Row 4 - 012474 / 大成成长回报六个月持有期混合型证券投资基金C
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "012474"
$newSheet.Range("C4").Value = "大成成长回报六个月持有期混合型证券投资基金C"
$newSheet.Range("D4").Value = "0.43"
$newSheet.Range("E4").Value = "71.30"
$newSheet.Range("F4").Value = "2.49"
$newSheet.Range("G4").Value = "0.0107"
$newSheet.Range("H4").Value = 10

# Row 5 - 011074 / 鹏华安润混合C
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "011074"
$newSheet.Range("C5").Value = "鹏华安润混合C"
$newSheet.Range("D5").Value = "0.25"
$newSheet.Range("E5").Value = "29.79"
$newSheet.Range("F5").Value = "1.28"
$newSheet.Range("G5").Value = "0.0032"
$newSheet.Range("H5").Value = 3

# Now that the sheet is fully populated, give it its real name and move it
# into place: right after "2021-Q4" / right before "总计".
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet.Move($totalSheet)

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row 2 for the 2022-Q1 totals,
#    pushing 2021-Q4 (was row 2) and 2021-Q3 (was row 3) down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert()

# Insert() copies the format of the row above (the header) onto the new
# row; clear that back off the plain data cells, then re-apply the
# correct "index column" style to A2 by copying it from an existing data
# row (A3, which already carries the right style after the shift).
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.28

# Re-sequence the index column for the rows that got pushed down (they
# kept their old 0/1 values; they need to become 1/2).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
